# Update run to cross multiple browsers: bump VERSION numbers and flip
# EXECUTE flags to "yes" on the TESTDATA sheet, then move the active
# selection from G5 to F5.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TESTDATA")

# VERSION column (C): 103.3 -> 109 / 110 / 111
$ws.Range("C2").Value = 109
$ws.Range("C3").Value = 110
$ws.Range("C4").Value = 111

# EXECUTE column (D): run all three browsers now
$ws.Range("D2").Value = "yes"
$ws.Range("D4").Value = "yes"

# Move the selection like it was left after the edit
$ws.Activate()
$ws.Range("F5").Select()
